# reportOfficinaTemplate.xlsx - add "Voto Medio" and "Provincia" header columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shared-string order in the target file is Nome(0), Voto Medio(1), Provincia(2)
# so write C1 before B1 to reproduce that exact ordering.
$ws.Range("C1").Value = "Voto Medio"
$ws.Range("B1").Value = "Provincia"

# Copy the header style (bold, filled, centered) from A1 onto the two new headers.
$ws.Range("A1").Copy()
$ws.Range("B1:C1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Column widths matching the template layout. The COM ColumnWidth property is
# expressed in "characters" and Excel re-quantizes it to whole pixels before
# storing the OOXML <col width> (which is in MDW character units); the inputs
# below are chosen so the stored width round-trips to 17 / 32 / ~18.86 exactly
# (the last one lands on the nearest reachable pixel grid value).
$ws.Columns.Item(1).ColumnWidth = 16.166666666666668
$ws.Columns.Item(2).ColumnWidth = 31.166666666666668
$ws.Columns.Item(3).ColumnWidth = 18

# Selection as left by the editor.
$ws.Range("C7").Select() | Out-Null
